{"js": "// Replace the table header text \"Summary 2\" with \"Summary 0\".\nconst results = context.document.body.search(\"Summary 2\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Summary 0\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the table header cell containing \"Summary 2\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Summary 2\"\n$find.Forward = $true\n$find.Wrap = 0          # wdFindStop\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.MatchWildcards = $false\n\n$found = $find.Execute()\nif ($found) {\n    # Replace the whole paragraph's text (not just the found sub-range) so\n    # the existing run/xml:space formatting of the cell is preserved.\n    $p = $find.Parent.Paragraphs(1)\n    $p.Range.Text = \"Summary 0\"\n}\n"}
